$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 31.444445
$ws.Range("I38").Value = 23
$ws.Range("J38").Value = 99
$ws.Range("K38").Value = 69
$ws.Range("L38").Value = 297
$ws.Range("M38").Value = 303
$ws.Range("N38").Value = -1041

$ws.Range("H40").Value = 2331.7144
$ws.Range("I40").Value = 2133.3333
$ws.Range("K40").Value = 2133.3333
$ws.Range("M40").Value = -1958.3333

$ws.Range("H58").Value = 106.666664
$ws.Range("I58").Value = 106.666664
$ws.Range("K58").Value = 319.999992
$ws.Range("M58").Value = -169.999992

$ws.Range("H107").Value = 75555.336
$ws.Range("I107").Value = 100385.89
$ws.Range("K107").Value = 100385.89
$ws.Range("M107").Value = -98465.89

$ws.Range("H138").Value = 5788.9697
$ws.Range("I138").Value = 4237.375
$ws.Range("J138").Value = 6285.48
$ws.Range("K138").Value = 12712.125
$ws.Range("L138").Value = 18856.44
$ws.Range("M138").Value = -7572.125
$ws.Range("N138").Value = -29136.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 140000
$ws.Range("J129").Value = 140000
$ws.Range("L129").Value = 140000
$ws.Range("N129").Value = -150000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -753

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H81").Value = 29166.166
$ws.Range("J81").Value = 29166.166
$ws.Range("L81").Value = 29166.166
$ws.Range("N81").Value = -31288.166

$ws.Range("H84").Value = 29166.166
$ws.Range("J84").Value = 29166.166
$ws.Range("L84").Value = 87498.49800000001
$ws.Range("N84").Value = -98106.49800000001

$ws.Range("H86").Value = 3333
$ws.Range("I86").Value = 3333
$ws.Range("K86").Value = 3333
$ws.Range("M86").Value = -2210

$ws.Range("H89").Value = 3333
$ws.Range("I89").Value = 3333
$ws.Range("K89").Value = 16665
$ws.Range("M89").Value = -11049

$ws.Range("H107").Value = 30705.572
$ws.Range("J107").Value = 2897.3333
$ws.Range("L107").Value = 2897.3333
$ws.Range("N107").Value = -6737.3333

$ws.Range("H117").Value = 49500
$ws.Range("J117").Value = 49500
$ws.Range("L117").Value = 49500
$ws.Range("N117").Value = -58678

$ws.Range("H135").Value = 39950
$ws.Range("J135").Value = 39950
$ws.Range("L135").Value = 39950
$ws.Range("N135").Value = -50090

$ws.Range("H138").Value = 70520
$ws.Range("J138").Value = 70520
$ws.Range("L138").Value = 70520
$ws.Range("N138").Value = -80800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H60").Value = 15497.5
$ws.Range("J60").Value = 13996.667
$ws.Range("L60").Value = 13996.667
$ws.Range("N60").Value = -15018.667

$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 2666.6667
$ws.Range("J62").Value = 1500
$ws.Range("L62").Value = 1500
$ws.Range("N62").Value = -2748

$ws.Range("H65").Value = 2666.6667
$ws.Range("J65").Value = 1500
$ws.Range("L65").Value = 7500
$ws.Range("N65").Value = -13740

$ws.Range("H74").Value = 41666.332
$ws.Range("J74").Value = 41749.5
$ws.Range("L74").Value = 41749.5
$ws.Range("N74").Value = -43497.5

$ws.Range("H77").Value = 41666.332
$ws.Range("J77").Value = 41749.5
$ws.Range("L77").Value = 125248.5
$ws.Range("N77").Value = -133984.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4302
$ws.Range("J5").Value = 4302
$ws.Range("L5").Value = 12906
$ws.Range("N5").Value = -13130

$ws.Range("H34").Value = 1843.1
$ws.Range("I34").Value = 108
$ws.Range("K34").Value = 324
$ws.Range("M34").Value = -240

$ws.Range("H39").Value = 1667.6666
$ws.Range("I39").Value = 1501.5
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 4504.5
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -4210.5
$ws.Range("N39").Value = -6588

$ws.Range("H55").Value = 3771.2239
$ws.Range("J55").Value = 3882.0317
$ws.Range("L55").Value = 11646.0951
$ws.Range("N55").Value = -12000.0951

$ws.Range("H113").Value = 847
$ws.Range("J113").Value = 808.75
$ws.Range("L113").Value = 2426.25
$ws.Range("N113").Value = -6766.25

$ws.Range("H135").Value = 4302
$ws.Range("J135").Value = 4302
$ws.Range("L135").Value = 38718
$ws.Range("N135").Value = -43788

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H102").Value = 1563.3334
$ws.Range("I102").Value = 1069.8572
$ws.Range("K102").Value = 1069.8572
$ws.Range("M102").Value = 552.1428000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 999
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("M46").Value = -811

$ws.Range("H61").Value = 2004
$ws.Range("I61").Value = 1927.7693
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 1927.7693
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -1725.7693
$ws.Range("N61").Value = -2903.5

$ws.Range("H113").Value = 2004
$ws.Range("I113").Value = 1927.7693
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 1927.7693
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = 242.2307000000001
$ws.Range("N113").Value = -6839.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996

$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -84984

$ws.Range("H96").Value = 5333.3335
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 4000
$ws.Range("N96").Value = -6746

$ws.Range("H113").Value = 374.33334
$ws.Range("I113").Value = 411.5
$ws.Range("K113").Value = 1234.5
$ws.Range("M113").Value = 935.5

$ws.Range("H132").Value = 400
$ws.Range("I132").Value = 400
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1200
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("M132").Value = 1330

$ws.Range("H136").Value = 3031.75
$ws.Range("I136").Value = 3031.75
$ws.Range("K136").Value = 9095.25
$ws.Range("M136").Value = -6545.25
